$d = $word.ActiveDocument

# 1. Insert a new bullet paragraph before "Change the way allocators are setup..."
#    with the same ListParagraph / numId=6 / ilvl=0 formatting, and give it its text.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Change the way allocators are setup*") {
        $targetIndex = $i
        break
    }
}
$d.Paragraphs($targetIndex).Range.InsertParagraphBefore()
$d.Paragraphs($targetIndex).Range.Text = "Look at destructor handling again, original design failed to account for user attempts to invoke the destructor of a non-GC allocated object."

# 2. Edit the "Minimize Impact" bullet: "final impact by" -> "final end-user code by"
$d.Content.Find.Execute("final impact by using strategies", $true, $false, $false, $false, $false, $true, 1, $false, "final end-user code by using strategies", 2) | Out-Null

# 3. Split the "Destructors will have 2 possible places..." run so that the
#    `_GoBack` bookmark sits right after "...will skip the call to `gc.markF"
$destructorsIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Destructors will have 2 possible places*") {
        $destructorsIndex = $i
        break
    }
}
$destructorsPara = $d.Paragraphs($destructorsIndex)
$splitPoint = $destructorsPara.Range.Start + 303
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
